# Apply "a bunch of updates" to the field_names_validated sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-categorize the various "Section 8" rows under the broader "Voucher" category.
$ws.Range("D32:D39").Value = "Voucher"

# "Shared Living (Not friends or relatives)" now counts toward Permanent Housing.
$ws.Range("E48").Value = "Permanent Housing"

# Clear out the stale reviewer notes/questions in column F that have since been
# resolved or are no longer relevant.
$ws.Range("F11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("F48").ClearContents()
$ws.Range("F54").ClearContents()
